# Added missing NOS0-UK00 transmission line to transferdata
#
# The "transferdata" sheet holds a sorted (by from-to) table of transmission
# lines, once per scenario-year (2030 block, then 2040 block). A row for the
# NOS0-UK00 connection was missing from both year blocks; this script inserts
# it in the correct alphabetically-sorted position in each block.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("transferdata")

# --- 2030 block: insert NOS0-UK00 right after NOS0-SE03 (row 42), before PL00-SE04 (row 43) ---
$ws.Rows.Item(43).Insert()
$ws.Range("A43").Value = "NOS0-UK00"
$ws.Range("B43").Value = "NOS0"
$ws.Range("C43").Value = "UK00"
$ws.Range("D43").Value = "elec"
$ws.Range("E43").Value = "National Trends"
$ws.Range("F43").Value = 2030
$ws.Range("G43").Value = 1400
$ws.Range("H43").Value = 1400
$ws.Range("I43").Value = 0.002768
$ws.Range("J43").Value = 1
$ws.Range("K43").Value = 0.01

# --- 2040 block: insert NOS0-UK00 right after NOS0-SE03 (row 88), before PL00-SE04 (row 89) ---
# (row numbers shifted down by 1 because of the insert above)
$ws.Rows.Item(89).Insert()
$ws.Range("A89").Value = "NOS0-UK00"
$ws.Range("B89").Value = "NOS0"
$ws.Range("C89").Value = "UK00"
$ws.Range("D89").Value = "elec"
$ws.Range("E89").Value = "National Trends"
$ws.Range("F89").Value = 2040
$ws.Range("G89").Value = 1400
$ws.Range("H89").Value = 1400
$ws.Range("I89").Value = 0.002768
$ws.Range("J89").Value = 1
$ws.Range("K89").Value = 0.01

# Grow the "Table6" structured table so it covers the two new rows (A1:K91 -> A1:K93)
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:K93"))

# --- restore/update per-sheet selections on the sheets whose cursor moved ---
$wsUnitdata = $wb.Worksheets.Item("unitdata")
$wsUnitdata.Range("M380").Select()

$wsUnittypedata = $wb.Worksheets.Item("unittypedata")
$wsUnittypedata.Range("A22").Select()

$wsEmissiondata = $wb.Worksheets.Item("emissiondata")
$wsEmissiondata.Range("K18").Select()

# Leave "transferdata" as the active sheet/tab, cursor on the newly-added 2040 row
$ws.Range("A89").Select()
$ws.Activate()
